$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(0, 300.5, 648.5218620877475)
    ,@(1, 314.5, 542.8372192339825)
    ,@(2, 328.4, 544.0243248082677)
    ,@(3, 342, 621.4941140871844)
    ,@(4, 356, 621.6624028950735)
    ,@(5, 369.6, 639.2826832849474)
    ,@(6, 383.5, 699.6753702954717)
    ,@(7, 397.5, 664.6434664923096)
    ,@(8, 411.1, 650.8052907531203)
    ,@(9, 425, 684.1310767041753)
    ,@(10, 439, 687.612527645598)
    ,@(11, 452.6, 675.3806862879186)
    ,@(12, 466.6, 682.5433297974471)
    ,@(13, 480.5, 687.4641506234847)
    ,@(14, 494.1, 708.2026844185129)
    ,@(15, 508.3, 723.5555709992348)
    ,@(16, 522, 743.1376355206964)
    ,@(17, 535.5999999999999, 737.9579287208893)
    ,@(18, 550, 740.0846604397326)
    ,@(19, 563.5, 747.0217090528308)
    ,@(20, 588.4000000000001, 758.0243611096105)
    ,@(21, 614, 769.913414766954)
    ,@(22, 639, 781.3338713095197)
    ,@(23, 664.6, 788.4694833929507)
    ,@(24, 689.5, 795.9990780785909)
    ,@(25, 715, 803.0749511449553)
    ,@(26, 740, 813.2477400347726)
    ,@(27, 765, 819.4318787552169)
    ,@(28, 790.5, 823.5963340752495)
    ,@(29, 816, 831.1378702148264)
    ,@(30, 841.3000000000001, 835.6419846562435)
    ,@(31, 866, 839.9682791802032)
    ,@(32, 891.5999999999999, 843.8533291053041)
    ,@(33, 917, 846.9174244863931)
    ,@(34, 942.4, 849.6832873156058)
    ,@(35, 967.3000000000001, 852.9014591023288)
    ,@(36, 992.7, 854.7578712015946)
    ,@(37, 1018, 857.3687578116162)
    ,@(38, 1055, 861.2125890975486)
    ,@(39, 1083, 864.1524254465738)
    ,@(40, 1112, 864.7332298017396)
    ,@(41, 1143, 866.9465102695284)
    ,@(42, 1174, 868.5250025152192)
    ,@(43, 1205, 869.8939773679193)
    ,@(44, 1237, 870.4474055173383)
    ,@(45, 1271, 871.5300530509618)
    ,@(46, 1306, 873.093703903538)
    ,@(47, 1341, 873.8747631843457)
    ,@(48, 1377, 875.5791450362204)
    ,@(49, 1414, 877.8691746745837)
    ,@(50, 1452, 880.4007064967534)
    ,@(51, 1492, 884.794222982821)
    ,@(52, 1532, 889.0068599922696)
    ,@(53, 1573, 895.6441404246451)
    ,@(54, 1616, 902.0479640041711)
    ,@(55, 1659, 905.0071596286014)
    ,@(56, 1704, 906.3696520759113)
    ,@(57, 1751, 907.978728982739)
    ,@(58, 2098, 913.890219148783)
    ,@(59, 2938, 925.5417562294977)
    ,@(60, 3060, 927.1123280009119)
    ,@(61, 3184, 928.6624337616622)
    ,@(62, 3307, 929.8786262196268)
    ,@(63, 3430, 931.0921673471029)
)

foreach ($row in $data) {
    $idx = [int]$row[0]
    $b = $row[1]
    $c = $row[2]
    $r = $idx + 2
    $ws.Cells.Item($r, 1).Value = $idx
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
}

# Copy the formatting (bold/border/center-top alignment style) from the
# existing last formatted row (A39) down onto the newly added A40:A65 cells.
$ws.Range("A39").Copy()
$ws.Range("A40:A65").PasteSpecial(-4122)
